$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @(
        "How do you add a new vectdraw font",
        "llama3.2:latest",
        "To add a new vector draw font, click on `"Choose Font`" in the Font dialog box. This will open the Font dialog box where you can change the font type, style, size, and color. You can then select a new vector draw font from the available options."
    ),
    @(
        "How do you define new curve pens?",
        "llama3.2:latest",
        "To define a new curve pen, follow these steps:`n1. On the Curve Settings tab, select Curve Settings.`n2. The Curve Settings sidebar will open.`n3. On the Display tab, open the Pens node.`n4. Select the `"Use Customized Pens`" double arrow.`n5. The Pen Definition dialog box will open.`n6. Select an unused number in the `"Select Pen`" box.`n7. Enter a Pen Name (e.g., DotDashDash).`n8. Select the Pen Type for both the display and printing pens.`nNote: A pen definition consists of up to five segments of dot-combinations, with up to ten pixels on and up to ten pixels off per segment."
    ),
    @(
        "What track types are there",
        "llama3.2:latest",
        "There are six track types in GEO:`n1. Depth`n2. Linear`n3. Logarithmic`n4. Lithology`n5. % Lithology`n6. Text`n7. Qualitative"
    ),
    @(
        "What track types are there",
        "llama3.2:latest",
        "There are six track types in GEO:`n1. Depth`n2. Linear`n3. Logarithmic`n4. Lithology`n5. % Lithology`n6. Text`n7. Qualitative"
    )
)

$startRow = 126
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
